$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the two new ones in order ---
$wsRandom = $wb.Worksheets.Item(1)
$wsRandom.Name = "Random"

$wsNFL = $wb.Worksheets.Add([System.Type]::Missing, $wsRandom)
$wsNFL.Name = "NFL"

$wsGames = $wb.Worksheets.Add([System.Type]::Missing, $wsNFL)
$wsGames.Name = "Games"

# --- NFL sheet data (column A down, then column B rows 2-4; header B1 filled in last) ---
$wsNFL.Range("A1").Value = "City"
$wsNFL.Range("A2").Value = "Buffalo"
$wsNFL.Range("B2").Value = "Bills"
$wsNFL.Range("A3").Value = "Miami"
$wsNFL.Range("B3").Value = "Dolphins"
$wsNFL.Range("A4").Value = "Denver"
$wsNFL.Range("B4").Value = "Broncos"
$wsNFL.Columns.Item(2).ColumnWidth = 16.333333333333336

# --- Games sheet data ---
$wsGames.Range("B1").Value = "Developer"
$wsGames.Range("A1").Value = "Title"
$wsGames.Range("A2").Value = "Fallout 76"
$wsGames.Range("B2").Value = "Bethesda"
$wsGames.Range("A3").Value = "Half-Life"
$wsGames.Range("B3").Value = "Valve"
$wsGames.Range("A4").Value = "Doom"
$wsGames.Range("B4").Value = "id"
$wsGames.Range("A5").Value = "Duke Nukem"
$wsGames.Range("B5").Value = "3D Realms"
$wsGames.Columns.Item(1).ColumnWidth = 13.833333333333332
$wsGames.Columns.Item(2).ColumnWidth = 13.0

# --- NFL header added last (matches shared-string ordering in the target) ---
$wsNFL.Range("B1").Value = "Team Name"

# --- Selections on each sheet ---
$wsRandom.Range("F28").Select()
$wsNFL.Range("F12").Select()
$wsGames.Range("E30").Select()

# --- Games is the active/selected tab ---
$wsGames.Activate()
